$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.520999999999999
$ws.Range("C9").Value = -10.85
$ws.Range("D12").Value = -7.375999999999999
$ws.Range("C18").Value = -12
$ws.Range("C20").Value = -12.354
$ws.Range("D26").Value = -8.122999999999999
$ws.Range("C27").Value = -13.005
$ws.Range("D27").Value = -8.369999999999999
$ws.Range("D29").Value = -7.292
$ws.Range("D37").Value = -7.943000000000001
$ws.Range("D38").Value = -7.688
$ws.Range("D51").Value = -8.391000000000002
$ws.Range("D55").Value = -7.918000000000001
$ws.Range("C69").Value = -11.269
$ws.Range("D69").Value = -7.243
$ws.Range("D70").Value = -7.203
$ws.Range("C76").Value = -13.032
$ws.Range("C82").Value = -11.991
$ws.Range("D83").Value = -8.373000000000001
$ws.Range("D102").Value = -7.861
